# Update "想去人数" (interested-count) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 126
$ws1.Range("F5").Value = 2823
$ws1.Range("F6").Value = 281
$ws1.Range("F7").Value = 391

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 126
$ws4.Range("F5").Value = 2823
$ws4.Range("F6").Value = 281
$ws4.Range("F9").Value = 391
